$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.682.45"
$ws.Range("E2").Value = "  -2.08%  "
$ws.Range("D3").Value = "2.339.17"
$ws.Range("E3").Value = "  -2.38%  "
$ws.Range("E4").Value = "  -0.01%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "502.85"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.45%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "128.55"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -2.97%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -2.96%  "
$ws.Range("D9").Value = "2.345.96"
$ws.Range("E9").Value = "  -2.32%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.0977"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("E12").Value = "  +3.46%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.318"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -1.68%  "
$ws.Range("D14").Value = "2.751.06"
$ws.Range("E14").Value = "  -2.42%  "
$ws.Range("D15").Value = "55.628.86"
$ws.Range("E15").Value = "  -2.04%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "21.54"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("E17").Value = "  -1.91%  "
$ws.Range("D18").Value = "2.365.31"
$ws.Range("E18").Value = "  -1.18%  "
$ws.Range("E19").Value = "  -2.83%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "309.76"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("E21").Value = "  -2.04%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.20"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.99%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.01%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "65.15"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -3.49%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("E26").Value = "  -1.18%  "
$ws.Range("E27").Value = "  -3.18%  "
$ws.Range("E28").Value = "  -4.59%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "171.46"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -2.30%  "
$ws.Range("E30").Value = "  -1.08%  "
$ws.Range("D31").Value = "0.0₃0701"
$ws.Range("E31").Value = "  -3.08%  "
$ws.Range("E33").Value = "  -1.46%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.997"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("E35").Value = "  -5.50%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "17.62"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -1.46%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.16"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -2.33%  "
$ws.Range("E38").Value = "  -4.57%  "
$ws.Range("E39").Value = "  -0.70%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "36.08"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -2.04%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.38"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -4.33%  "
$ws.Range("E42").Value = "  -0.77%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "126.10"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -3.97%  "
$ws.Range("E44").Value = "  -3.55%  "
$ws.Range("E45").Value = "  -2.63%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0891"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -2.31%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "236.54"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -5.65%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.0473"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -2.89%  "
$ws.Range("E49").Value = "  -2.71%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "16.77"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.56%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.952"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.13%  "
